# Updates the FFXIV Leve profit calculations across all class sheets
# (prices/profits recalculated from refreshed market data).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2187.6428
$ws.Range("I19").Value = 1161.2858
$ws.Range("J19").Value = 3214
$ws.Range("K19").Value = 1161.2858
$ws.Range("L19").Value = 3214
$ws.Range("M19").Value = -986.2858000000001
$ws.Range("N19").Value = -3564
$ws.Range("H116").Value = 5085
$ws.Range("I116").Value = 5127.5
$ws.Range("K116").Value = 5127.5
$ws.Range("M116").Value = -1685.5
$ws.Range("H138").Value = 9261147
$ws.Range("I138").Value = 1377.2858
$ws.Range("J138").Value = 15153728
$ws.Range("K138").Value = 4131.857400000001
$ws.Range("L138").Value = 45461184
$ws.Range("M138").Value = 1008.142599999999
$ws.Range("N138").Value = -45471464

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2294.8572
$ws.Range("I45").Value = 2007.2
$ws.Range("J45").Value = 3014
$ws.Range("K45").Value = 2007.2
$ws.Range("L45").Value = 3014
$ws.Range("M45").Value = -1630.2
$ws.Range("N45").Value = -3768
$ws.Range("H61").Value = 47624052
$ws.Range("I61").Value = 66670304
$ws.Range("K61").Value = 66670304
$ws.Range("M61").Value = -66670092
$ws.Range("H63").Value = 2691.5435
$ws.Range("I63").Value = 1465.1765
$ws.Range("K63").Value = 1465.1765
$ws.Range("M63").Value = -779.1765
$ws.Range("H66").Value = 2691.5435
$ws.Range("I66").Value = 1465.1765
$ws.Range("K66").Value = 7325.8825
$ws.Range("M66").Value = -3893.8825
$ws.Range("H132").Value = 43482080
$ws.Range("I132").Value = 4041.1
$ws.Range("K132").Value = 12123.3
$ws.Range("M132").Value = -9593.299999999999
$ws.Range("H135").Value = 31746
$ws.Range("J135").Value = 31746
$ws.Range("L135").Value = 31746
$ws.Range("N135").Value = -41886
$ws.Range("H136").Value = 47624052
$ws.Range("I136").Value = 66670304
$ws.Range("K136").Value = 200010912
$ws.Range("M136").Value = -200008362

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 13882.25
$ws.Range("I105").Value = 17809.666
$ws.Range("K105").Value = 17809.666
$ws.Range("M105").Value = -16062.666
$ws.Range("H134").Value = 3005.425
$ws.Range("I134").Value = 2873.4595
$ws.Range("J134").Value = 4633
$ws.Range("K134").Value = 8620.378499999999
$ws.Range("L134").Value = 13899
$ws.Range("M134").Value = -6085.378499999999
$ws.Range("N134").Value = -18969

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 6000
$ws.Range("I11").Value = 6000
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 6000
$ws.Range("L11").Value = 0
$ws.Range("N11").Value = $null
$ws.Range("M11").Value = -5860
$ws.Range("H31").Value = 23260198
$ws.Range("I31").Value = 3348.3928
$ws.Range("J31").Value = 66672984
$ws.Range("K31").Value = 3348.3928
$ws.Range("L31").Value = 66672984
$ws.Range("M31").Value = -3053.3928
$ws.Range("N31").Value = -66673574
$ws.Range("H34").Value = 23260198
$ws.Range("I34").Value = 3348.3928
$ws.Range("J34").Value = 66672984
$ws.Range("K34").Value = 3348.3928
$ws.Range("L34").Value = 66672984
$ws.Range("M34").Value = -3146.3928
$ws.Range("N34").Value = -66673388
$ws.Range("H52").Value = 109750.6
$ws.Range("J52").Value = 112189
$ws.Range("L52").Value = 112189
$ws.Range("N52").Value = -112777
$ws.Range("H94").Value = 1869.3572
$ws.Range("I94").Value = 1568.25
$ws.Range("J94").Value = 1989.8
$ws.Range("K94").Value = 1568.25
$ws.Range("L94").Value = 1989.8
$ws.Range("M94").Value = -1117.25
$ws.Range("N94").Value = -2891.8
$ws.Range("H105").Value = 10964.5
$ws.Range("I105").Value = 2460.5
$ws.Range("K105").Value = 2460.5
$ws.Range("M105").Value = -713.5
$ws.Range("H107").Value = 1386.6428
$ws.Range("I107").Value = 801.8
$ws.Range("K107").Value = 801.8
$ws.Range("M107").Value = 1118.2
$ws.Range("H122").Value = 1664.0358
$ws.Range("I122").Value = 1535.5714
$ws.Range("K122").Value = 4606.7142
$ws.Range("M122").Value = -2156.7142
$ws.Range("H134").Value = 1233.3
$ws.Range("I134").Value = 1037.1666
$ws.Range("K134").Value = 3111.4998
$ws.Range("M134").Value = -576.4998000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 21016.084
$ws.Range("I56").Value = 21016.084
$ws.Range("K56").Value = 21016.084
$ws.Range("M56").Value = -20486.084
$ws.Range("H86").Value = 1041.2858
$ws.Range("I86").Value = 998.5
$ws.Range("J86").Value = 1058.4
$ws.Range("K86").Value = 2995.5
$ws.Range("L86").Value = 3175.2
$ws.Range("M86").Value = -1809.5
$ws.Range("N86").Value = -5547.200000000001
$ws.Range("H87").Value = 1273.2
$ws.Range("J87").Value = 4016
$ws.Range("L87").Value = 12048
$ws.Range("N87").Value = -14544
$ws.Range("H89").Value = 1041.2858
$ws.Range("I89").Value = 998.5
$ws.Range("J89").Value = 1058.4
$ws.Range("K89").Value = 8986.5
$ws.Range("L89").Value = 9525.6
$ws.Range("M89").Value = -3058.5
$ws.Range("N89").Value = -21381.6
$ws.Range("H90").Value = 1273.2
$ws.Range("J90").Value = 4016
$ws.Range("L90").Value = 36144
$ws.Range("N90").Value = -48624
$ws.Range("H113").Value = 3896.3635
$ws.Range("J113").Value = 5310.1665
$ws.Range("L113").Value = 15930.4995
$ws.Range("N113").Value = -20270.4995
$ws.Range("H131").Value = 1760.5454
$ws.Range("I131").Value = 1402.1
$ws.Range("J131").Value = 1916.3914
$ws.Range("K131").Value = 4206.299999999999
$ws.Range("L131").Value = 5749.174199999999
$ws.Range("M131").Value = 833.7000000000007
$ws.Range("N131").Value = -15829.1742

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3435.261
$ws.Range("I113").Value = 2573.9333
$ws.Range("K113").Value = 2573.9333
$ws.Range("M113").Value = -403.9333000000001
$ws.Range("H122").Value = 7772.4375
$ws.Range("I122").Value = 3895.6
$ws.Range("K122").Value = 11686.8
$ws.Range("M122").Value = -9236.799999999999
$ws.Range("H126").Value = 25007240
$ws.Range("H132").Value = 3743.6667
$ws.Range("I132").Value = 3586.2307
$ws.Range("K132").Value = 10758.6921
$ws.Range("M132").Value = -8228.6921

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H50").Value = 35999
$ws.Range("I50").Value = 40999
$ws.Range("K50").Value = 40999
$ws.Range("M50").Value = -40362
$ws.Range("H68").Value = 2864.4285
$ws.Range("I68").Value = 2633.5833
$ws.Range("J68").Value = 4249.5
$ws.Range("K68").Value = 2633.5833
$ws.Range("L68").Value = 4249.5
$ws.Range("M68").Value = -1884.5833
$ws.Range("N68").Value = -5747.5
$ws.Range("H71").Value = 2864.4285
$ws.Range("I71").Value = 2633.5833
$ws.Range("J71").Value = 4249.5
$ws.Range("K71").Value = 13167.9165
$ws.Range("L71").Value = 21247.5
$ws.Range("M71").Value = -9423.916499999999
$ws.Range("N71").Value = -28735.5
$ws.Range("H82").Value = 4959.3
$ws.Range("I82").Value = 2933.1667
$ws.Range("J82").Value = 7998.5
$ws.Range("K82").Value = 2933.1667
$ws.Range("L82").Value = 7998.5
$ws.Range("M82").Value = -2572.1667
$ws.Range("N82").Value = -8720.5
$ws.Range("H85").Value = 4959.3
$ws.Range("I85").Value = 2933.1667
$ws.Range("J85").Value = 7998.5
$ws.Range("K85").Value = 2933.1667
$ws.Range("L85").Value = 7998.5
$ws.Range("M85").Value = -1685.1667
$ws.Range("N85").Value = -10494.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7122.3335
$ws.Range("J62").Value = 7179.615
$ws.Range("L62").Value = 7179.615
$ws.Range("N62").Value = -8427.615
$ws.Range("H65").Value = 7122.3335
$ws.Range("J65").Value = 7179.615
$ws.Range("L65").Value = 35898.075
$ws.Range("N65").Value = -42138.075
$ws.Range("H96").Value = 7642.1665
$ws.Range("I96").Value = 5980
$ws.Range("J96").Value = 8829.429
$ws.Range("K96").Value = 5980
$ws.Range("L96").Value = 8829.429
$ws.Range("M96").Value = -4607
$ws.Range("N96").Value = -11575.429
$ws.Range("H113").Value = 665.2273
$ws.Range("I113").Value = 261.5625
$ws.Range("K113").Value = 784.6875
$ws.Range("M113").Value = 1385.3125
$ws.Range("H122").Value = 52633536
$ws.Range("I122").Value = 62501636
$ws.Range("K122").Value = 187504908
$ws.Range("M122").Value = -187502458
$ws.Range("H123").Value = 57882.668
$ws.Range("J123").Value = 57882.668
$ws.Range("L123").Value = 57882.668
$ws.Range("N123").Value = -67682.66800000001
$ws.Range("H132").Value = 4458.857
$ws.Range("I132").Value = 4458.857
$ws.Range("K132").Value = 13376.571
$ws.Range("M132").Value = -10846.571
